# date overlap logic added
# Swap "At Work" (C) / "Annual Leave" (G) marks for 18-20 and 30-33,
# and add "Sick Leave" (E) marks / remove "Annual Leave" (G) marks for 23-27.
# Update the Total row (44) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18-20: move mark from G (Annual Leave) to C (At Work)
foreach ($r in 18..20) {
    $ws.Range("C$r").Value = 1
    $ws.Range("G$r").Value = $null
}

# Rows 23-27: move mark from G (Annual Leave) to E (Sick Leave)
foreach ($r in 23..27) {
    $ws.Range("E$r").Value = 1
    $ws.Range("G$r").Value = $null
}

# Rows 30-33: move mark from C (At Work) to G (Annual Leave)
foreach ($r in 30..33) {
    $ws.Range("C$r").Value = $null
    $ws.Range("G$r").Value = 1
}

# Update Total row 44
$ws.Range("C44").Value = 11
$ws.Range("E44").Value = 5
$ws.Range("G44").Value = 4
